$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Name" / value row at the top of the form (row 2)
$ws.Range("A2").Value = "Name"
$ws.Range("B2").Value = "SPREADSHEETFORM:SINGLE:name/value"

# Move the active selection to B3 (as in the edited workbook)
$ws.Range("B3").Select()
